$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$checkStyle = $ws.Range("K3").Style

$rows = @(19, 30, 31, 32, 33, 35, 37)
foreach ($r in $rows) {
    $cell = $ws.Range("K$r")
    $cell.Value = "✓"
    $cell.Style = $checkStyle
}

$ws.Range("K19").Select()
